$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "['MCT-3A-Elementos de máquinas', -, -, -]"

# Row 3
$ws.Range("E3").Value = "['MCT-3A-Elementos de máquinas', -, -, -]"

# Row 4
$ws.Range("E4").Value = "['MCT-3A-Elementos de máquinas', -, -, -]"

# Row 6
$ws.Range("B6").Value = "-"
$ws.Range("E6").Value = "['MCT-3A-Elementos de máquinas', -, -, -]"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("E7").Value = "MCT-2A-M.T.R.M."

# Row 8
$ws.Range("E8").Value = "MCT-2A-M.T.R.M."
